# Generate and save output file after processing
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns for general_college_subjects.history/electives/cs
# right before the existing general_college_subjects.arts column (R),
# shifting everything from R onward to the right.
$ws.Range("R1:T1").Insert(-4161)

# New header labels for the inserted columns
$ws.Range("R1").Value = "general_college_subjects.history"
$ws.Range("S1").Value = "general_college_subjects.electives"
$ws.Range("T1").Value = "general_college_subjects.cs"

# New data values for row 2 in the inserted columns
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0

# Update previously "Unknown" placeholder values with the processed results
$ws.Range("D2").Value = "not considered"
$ws.Range("E2").Value = "considered"
$ws.Range("F2").Value = "considered"
$ws.Range("G2").Value = "considered"
$ws.Range("H2").Value = "very important"
$ws.Range("I2").Value = "considered"
$ws.Range("J2").Value = "considered"
